$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 21 de Mayo de 2020 a las 17:05"

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 1597124
$ws.Range("C4").Value = 4401
$ws.Range("D4").Value = 371077
$ws.Range("E4").Value = 1130935
$ws.Range("G4").Value = 176
$ws.Range("H4").Value = 95112

# Row 14 - India
$ws.Range("B14").Value = 114478
$ws.Range("C14").Value = 2450
$ws.Range("E14").Value = 65011
$ws.Range("G14").Value = 31
$ws.Range("H14").Value = 3465

# Row 17 - Canada
$ws.Range("B17").Value = 80555
$ws.Range("C17").Value = 413
$ws.Range("E17").Value = 33717
$ws.Range("G17").Value = 31
$ws.Range("H17").Value = 6062

# Row 30 - Singapur
$ws.Range("D30").Value = 12117
$ws.Range("E30").Value = 17672
$ws.Range("G30").Value = 1
$ws.Range("H30").Value = 23

# Row 133 - Reunion
$ws.Range("B133").Value = 449
$ws.Range("C133").Value = 2
$ws.Range("E133").Value = 37

# Row 156 - Mozambique
$ws.Range("B156").Value = 162
$ws.Range("C156").Value = 6
$ws.Range("E156").Value = 114

# Row 166 - Camboya
$ws.Range("B166").Value = 123
$ws.Range("C166").Value = 1
$ws.Range("E166").Value = 1

# Row 176 - Libia
$ws.Range("B176").Value = 71
$ws.Range("C176").Value = 2
$ws.Range("E176").Value = 33
